# Updated symbol list on Fri Dec 30 14:30:49 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking price as literal TEXT (matches the sheet's
# existing convention where Price/Volume columns are stored as text, not numbers).
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
}

# --- Rows 4-17: coin list shifted up by one (LEO moved from row 4 to row 17) ---
$rowUpdates = @{
    4 = @{ B = "HuobiToken"; C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D = "5.176"; E = "3HuobiTokenHT" }
    5 = @{ B = "Cronos"; C = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D = "0.05698"; E = "4CronosCRO" }
    6 = @{ B = "KuCoinToken"; C = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"; D = "6.485"; E = "5KuCoinTokenKCS" }
    7 = @{ B = "GateToken"; C = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"; D = "2.955"; E = "6GateTokenGT" }
    8 = @{ B = "MXToken"; C = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D = "0.8098"; E = "7MXTokenMX" }
    9 = @{ B = "FTXToken"; C = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"; D = "0.8351"; E = "8FTXTokenFTT" }
    10 = @{ B = "WazirX"; C = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"; D = "0.1336"; E = "9WazirXWRX" }
    11 = @{ B = "MandalaExchangeToken"; C = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"; D = "0.06954"; E = "10MandalaExchangeTokenMDX" }
    12 = @{ B = "BitrueCoin"; C = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"; D = "0.02827"; E = "11BitrueCoinBTR" }
    13 = @{ B = "BitMartToken"; C = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"; D = "0.09376"; E = "12BitMartTokenBMX" }
    14 = @{ B = "BitForexToken"; C = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"; D = "0.001511"; E = "13BitForexTokenBF" }
    15 = @{ B = "One"; C = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"; D = "0.0005950"; E = "14OneONE" }
    16 = @{ B = "TigerCash"; C = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"; D = "0.006188"; E = "15TigerCashTCH" }
    17 = @{ B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "3.500"; E = "16LEOLEO" }
}
foreach ($r in $rowUpdates.Keys) {
    $data = $rowUpdates[$r]
    $ws.Range("B$r").Value = $data.B
    $ws.Range("C$r").Value = $data.C
    Set-TextValue $ws.Range("D$r") $data.D
    $ws.Range("E$r").Value = $data.E
}

# --- Price-only (column D) updates ---
$priceUpdates = @{
    2 = "243.53"
    19 = "0.3194"
    20 = "0.03198"
    21 = "0.1337"
    22 = "3.750"
    23 = "0.04680"
    26 = "0.004244"
    27 = "0.00009701"
    40 = "0.03628"
    42 = "0.1049"
    44 = "0.007368"
    45 = "0.00005269"
    48 = "0.002296"
}
foreach ($r in $priceUpdates.Keys) {
    Set-TextValue $ws.Range("D$r") $priceUpdates[$r]
}

# --- Row 41: price + volume(1h) label update ---
Set-TextValue $ws.Range("D41") "0.006260"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"

# --- Row 43: volume(1h) label update ---
$ws.Range("E43").Value = "42CEJICEJI"

# --- Row 47: volume(1h) label update ---
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

